# Update the "Last Updated" timestamp columns across all sheets of the
# combined arrivals workbook from "2026-02-18 00:21" to "2026-02-18 00:28".

$wb = $excel.ActiveWorkbook

$oldValue = "2026-02-18 00:21"
$newValue = "2026-02-18 00:28"

# Sheet "Sydney": "Last Updated" values live in column H, rows 2-22.
$wsSydney = $wb.Worksheets.Item("Sydney")
for ($r = 2; $r -le 22; $r++) {
    $cell = $wsSydney.Cells.Item($r, 8)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}

# Sheet "Melbourne": "Last Updated" values live in column G, rows 2-84.
$wsMelbourne = $wb.Worksheets.Item("Melbourne")
for ($r = 2; $r -le 84; $r++) {
    $cell = $wsMelbourne.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}

# Sheet "Last Updated": single summary timestamp in B1.
$wsLastUpdated = $wb.Worksheets.Item("Last Updated")
$wsLastUpdated.Range("B1").Value = $newValue
